$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 6062582.5
$ws.Range("I33").Value = 3012.625
$ws.Range("J33").Value = 12987805
$ws.Range("K33").Value = 3012.625
$ws.Range("L33").Value = 12987805
$ws.Range("M33").Value = -2783.625
$ws.Range("N33").Value = -12988263
$ws.Range("H64").Value = 6100
$ws.Range("I64").Value = 7783.3335
$ws.Range("J64").Value = 3575
$ws.Range("K64").Value = 7783.3335
$ws.Range("L64").Value = 3575
$ws.Range("M64").Value = -7535.3335
$ws.Range("N64").Value = -4071
$ws.Range("H67").Value = 6100
$ws.Range("I67").Value = 7783.3335
$ws.Range("J67").Value = 3575
$ws.Range("K67").Value = 7783.3335
$ws.Range("L67").Value = 3575
$ws.Range("M67").Value = -6925.3335
$ws.Range("N67").Value = -5291
$ws.Range("H100").Value = 12821475
$ws.Range("I100").Value = 12821475
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 12821475
$ws.Range("L100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("N100").Value = -12820934
$ws.Range("H105").Value = 29333.334
$ws.Range("J105").Value = 29333.334
$ws.Range("L105").Value = 29333.334
$ws.Range("N105").Value = -36321.334
$ws.Range("H107").Value = 20835656
$ws.Range("I107").Value = 62501228
$ws.Range("J107").Value = 2872.125
$ws.Range("K107").Value = 62501228
$ws.Range("L107").Value = 2872.125
$ws.Range("M107").Value = -62499308
$ws.Range("N107").Value = -6712.125
$ws.Range("H115").Value = 380
$ws.Range("I115").Value = 380
$ws.Range("K115").Value = 1140
$ws.Range("M115").Value = 427
$ws.Range("H132").Value = 4168415.8
$ws.Range("I132").Value = 1320.6301
$ws.Range("K132").Value = 3961.8903
$ws.Range("M132").Value = -1431.8903
$ws.Range("H137").Value = 1194.5349
$ws.Range("I137").Value = 1234.3928
$ws.Range("J137").Value = 1120.1333
$ws.Range("K137").Value = 3703.1784
$ws.Range("L137").Value = 3360.3999
$ws.Range("M137").Value = -1153.1784
$ws.Range("N137").Value = -8460.3999
$ws.Range("H138").Value = 3413.6235
$ws.Range("I138").Value = 1661.8485
$ws.Range("J138").Value = 4377.1
$ws.Range("K138").Value = 4985.5455
$ws.Range("L138").Value = 13131.3
$ws.Range("M138").Value = 154.4544999999998
$ws.Range("N138").Value = -23411.3
$ws.Range("H139").Value = 51648.57
$ws.Range("J139").Value = 51648.57
$ws.Range("L139").Value = 51648.57
$ws.Range("N139").Value = -61928.57

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H27").Value = 14105
$ws.Range("J27").Value = 14105
$ws.Range("L27").Value = 14105
$ws.Range("N27").Value = -14473
$ws.Range("H32").Value = 17811.96
$ws.Range("I32").Value = 13943.362
$ws.Range("J32").Value = 31010.705
$ws.Range("K32").Value = 13943.362
$ws.Range("L32").Value = 31010.705
$ws.Range("M32").Value = -13656.362
$ws.Range("N32").Value = -31584.705
$ws.Range("H133").Value = 41440
$ws.Range("J133").Value = 41440
$ws.Range("L133").Value = 41440
$ws.Range("N133").Value = -46500

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 22013.531
$ws.Range("I134").Value = 1555.0454
$ws.Range("K134").Value = 4665.1362
$ws.Range("M134").Value = -2130.1362
$ws.Range("H135").Value = 64325.383
$ws.Range("J135").Value = 64325.383
$ws.Range("L135").Value = 64325.383
$ws.Range("N135").Value = -74465.383

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1278.826
$ws.Range("I107").Value = 870.6923
$ws.Range("J107").Value = 1809.4
$ws.Range("K107").Value = 870.6923
$ws.Range("L107").Value = 1809.4
$ws.Range("M107").Value = 1049.3077
$ws.Range("N107").Value = -5649.4
$ws.Range("H116").Value = 35000
$ws.Range("J116").Value = 35000
$ws.Range("L116").Value = 35000
$ws.Range("N116").Value = -44178
$ws.Range("H133").Value = 56659.6
$ws.Range("J133").Value = 56659.6
$ws.Range("L133").Value = 56659.6
$ws.Range("N133").Value = -61719.6
$ws.Range("H135").Value = 39028.26
$ws.Range("J135").Value = 39028.26
$ws.Range("L135").Value = 39028.26
$ws.Range("N135").Value = -49168.26

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1852879.8
$ws.Range("I131").Value = 6250684
$ws.Range("J131").Value = 1172.8158
$ws.Range("K131").Value = 18752052
$ws.Range("L131").Value = 3518.4474
$ws.Range("M131").Value = -18747012
$ws.Range("N131").Value = -13598.4474
$ws.Range("H132").Value = 4275779.5
$ws.Range("J132").Value = 7410332
$ws.Range("L132").Value = 66692988
$ws.Range("N132").Value = -66698048
$ws.Range("H137").Value = 16535.615
$ws.Range("J137").Value = 22605.076
$ws.Range("L137").Value = 67815.228
$ws.Range("N137").Value = -78015.228

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 427.36365
$ws.Range("I16").Value = 433.22223
$ws.Range("J16").Value = 401
$ws.Range("K16").Value = 433.22223
$ws.Range("L16").Value = 401
$ws.Range("M16").Value = -263.22223
$ws.Range("N16").Value = -741
$ws.Range("H46").Value = 1166.75
$ws.Range("I46").Value = 771.5714
$ws.Range("J46").Value = 1720
$ws.Range("K46").Value = 771.5714
$ws.Range("L46").Value = 1720
$ws.Range("M46").Value = -583.5714
$ws.Range("N46").Value = -2096
$ws.Range("H136").Value = 6241.8
$ws.Range("I136").Value = 2505.276
$ws.Range("J136").Value = 16092.637
$ws.Range("K136").Value = 7515.828
$ws.Range("L136").Value = 48277.911
$ws.Range("M136").Value = -4965.828
$ws.Range("N136").Value = -53377.911

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H99").Value = 25000
$ws.Range("J99").Value = 25000
$ws.Range("L99").Value = 25000
$ws.Range("N99").Value = -30990
$ws.Range("H117").Value = 26469.666
$ws.Range("J117").Value = 26469.666
$ws.Range("L117").Value = 26469.666
$ws.Range("N117").Value = -35647.666
$ws.Range("H126").Value = 934.26086
$ws.Range("I126").Value = 900.4706
$ws.Range("J126").Value = 1030
$ws.Range("K126").Value = 2701.4118
$ws.Range("L126").Value = 3090
$ws.Range("M126").Value = -231.4117999999999
$ws.Range("N126").Value = -8030
$ws.Range("H132").Value = 1690.1459
$ws.Range("I132").Value = 853.4138
$ws.Range("J132").Value = 2967.2632
$ws.Range("K132").Value = 2560.2414
$ws.Range("L132").Value = 8901.7896
$ws.Range("M132").Value = -30.24139999999989
$ws.Range("N132").Value = -13961.7896
$ws.Range("H136").Value = 2926047.2
$ws.Range("I136").Value = 1882.5238
$ws.Range("K136").Value = 5647.5714
$ws.Range("M136").Value = -3097.5714
